$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 2219.2727
$ws.Range("J112").Value = 2402.923
$ws.Range("L112").Value = 7208.768999999999
$ws.Range("N112").Value = -9424.769

# Row 135
$ws.Range("H135").Value = 3877.4285
$ws.Range("I135").Value = 1451.3158
$ws.Range("K135").Value = 13061.8422
$ws.Range("M135").Value = -10526.8422

# Row 137
$ws.Range("H137").Value = 4637.478
$ws.Range("I137").Value = 1305.64
$ws.Range("K137").Value = 3916.92
$ws.Range("M137").Value = -1366.92

# Row 138
$ws.Range("H138").Value = 2292.104
$ws.Range("I138").Value = 1853.85
$ws.Range("J138").Value = 2445.8772
$ws.Range("K138").Value = 5561.549999999999
$ws.Range("L138").Value = 7337.6316
$ws.Range("M138").Value = -421.5499999999993
$ws.Range("N138").Value = -17617.6316

# Row 141
$ws.Range("H141").Value = 3961.625
$ws.Range("I141").Value = 3248.625
$ws.Range("J141").Value = 4674.625
$ws.Range("K141").Value = 9745.875
$ws.Range("L141").Value = 14023.875
$ws.Range("M141").Value = -4565.875
$ws.Range("N141").Value = -24383.875

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 2666.6667
$ws.Range("I8").Value = 2000
$ws.Range("K8").Value = 2000
$ws.Range("M8").Value = -1856

# Row 11
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 32
$ws.Range("H32").Value = 166749.39
$ws.Range("I32").Value = 183240.39
$ws.Range("K32").Value = 183240.39
$ws.Range("M32").Value = -182953.39

# Row 45
$ws.Range("H45").Value = 4283.3335
$ws.Range("J45").Value = 5450
$ws.Range("L45").Value = 5450
$ws.Range("N45").Value = -6204

# Row 94
$ws.Range("H94").Value = 39544.625
$ws.Range("J94").Value = 39544.625
$ws.Range("L94").Value = 39544.625
$ws.Range("N94").Value = -41346.625

# Row 122
$ws.Range("H122").Value = 1714.5
$ws.Range("I122").Value = 1274.3334
$ws.Range("J122").Value = 2374.75
$ws.Range("K122").Value = 3823.0002
$ws.Range("L122").Value = 7124.25
$ws.Range("M122").Value = -1373.0002
$ws.Range("N122").Value = -12024.25

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 2591347.8
$ws.Range("J7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2226

# Row 99
$ws.Range("H99").Value = 17999.834
$ws.Range("I99").Value = 26125
$ws.Range("K99").Value = 26125
$ws.Range("M99").Value = -24627

# Row 107
$ws.Range("H107").Value = 822.8333
$ws.Range("I107").Value = 787.25
$ws.Range("K107").Value = 787.25
$ws.Range("M107").Value = 1132.75

$ws = $wb.Worksheets.Item("CRP")
# Row 11
$ws.Range("H11").Value = 330.16666
$ws.Range("I11").Value = 395
$ws.Range("K11").Value = 395
$ws.Range("M11").Value = -255

# Row 22
$ws.Range("H22").Value = 2056.353
$ws.Range("J22").Value = 2999
$ws.Range("L22").Value = 2999
$ws.Range("N22").Value = -3699

# Row 31
$ws.Range("H31").Value = 2334.8635
$ws.Range("I31").Value = 2702.25
$ws.Range("J31").Value = 1894
$ws.Range("K31").Value = 2702.25
$ws.Range("L31").Value = 1894
$ws.Range("M31").Value = -2407.25
$ws.Range("N31").Value = -2484

# Row 34
$ws.Range("H34").Value = 2334.8635
$ws.Range("I34").Value = 2702.25
$ws.Range("J34").Value = 1894
$ws.Range("K34").Value = 2702.25
$ws.Range("L34").Value = 1894
$ws.Range("M34").Value = -2500.25
$ws.Range("N34").Value = -2298

# Row 122
$ws.Range("H122").Value = 7009.727
$ws.Range("I122").Value = 1654.5358
$ws.Range("J122").Value = 36998.8
$ws.Range("K122").Value = 4963.607400000001
$ws.Range("L122").Value = 110996.4
$ws.Range("M122").Value = -2513.607400000001
$ws.Range("N122").Value = -115896.4

# Row 132
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7970
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 9385.84
$ws.Range("J114").Value = 10923.667
$ws.Range("L114").Value = 32771.001
$ws.Range("N114").Value = -39279.001

# Row 117
$ws.Range("H117").Value = 2575
$ws.Range("J117").Value = 2749.25
$ws.Range("L117").Value = 8247.75
$ws.Range("N117").Value = -15131.75

# Row 119
$ws.Range("H119").Value = 10789.75
$ws.Range("I119").Value = 3287.8333
$ws.Range("K119").Value = 9863.499899999999
$ws.Range("M119").Value = -5025.499899999999

# Row 132
$ws.Range("H132").Value = 1694
$ws.Range("I132").Value = 1394.75
$ws.Range("J132").Value = 1865
$ws.Range("K132").Value = 12552.75
$ws.Range("L132").Value = 16785
$ws.Range("M132").Value = -10022.75
$ws.Range("N132").Value = -21845

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2919.3333
$ws.Range("I102").Value = 2992.1428
$ws.Range("K102").Value = 2992.1428
$ws.Range("M102").Value = -1370.1428

# Row 113
$ws.Range("H113").Value = 4486.3
$ws.Range("I113").Value = 3459
$ws.Range("K113").Value = 3459
$ws.Range("M113").Value = -1289

# Row 126
$ws.Range("H126").Value = 16443.8
$ws.Range("I126").Value = 25073
$ws.Range("K126").Value = 75219
$ws.Range("M126").Value = -72749

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6922.7
$ws.Range("I7").Value = 4746.857
$ws.Range("K7").Value = 4746.857
$ws.Range("M7").Value = -4634.857

# Row 9
$ws.Range("H9").Value = 2228.75
$ws.Range("I9").Value = 466
$ws.Range("J9").Value = 5166.6665
$ws.Range("K9").Value = 466
$ws.Range("L9").Value = 5166.6665
$ws.Range("M9").Value = -242
$ws.Range("N9").Value = -5614.6665

# Row 22
$ws.Range("H22").Value = 2198.4167
$ws.Range("J22").Value = 2226.9714
$ws.Range("L22").Value = 2226.9714
$ws.Range("N22").Value = -2816.9714

# Row 27
$ws.Range("H27").Value = 2198.4167
$ws.Range("J27").Value = 2226.9714
$ws.Range("L27").Value = 2226.9714
$ws.Range("N27").Value = -2440.9714

# Row 126
$ws.Range("H126").Value = 6922.7
$ws.Range("I126").Value = 4746.857
$ws.Range("K126").Value = 14240.571
$ws.Range("M126").Value = -11770.571

# Row 132
$ws.Range("H132").Value = 2944.182
$ws.Range("I132").Value = 2222
$ws.Range("J132").Value = 3016.4
$ws.Range("K132").Value = 6666
$ws.Range("L132").Value = 9049.200000000001
$ws.Range("M132").Value = -4136
$ws.Range("N132").Value = -14109.2

# Row 136
$ws.Range("H136").Value = 7361.5713
$ws.Range("I136").Value = 2720.5715
$ws.Range("J136").Value = 9682.071
$ws.Range("K136").Value = 8161.7145
$ws.Range("L136").Value = 29046.213
$ws.Range("M136").Value = -5611.7145
$ws.Range("N136").Value = -34146.213

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 7500
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# Row 82
$ws.Range("H82").Value = 4750
$ws.Range("J82").Value = 4750
$ws.Range("L82").Value = 4750
$ws.Range("N82").Value = -5516

# Row 85
$ws.Range("H85").Value = 4750
$ws.Range("J85").Value = 4750
$ws.Range("L85").Value = 4750
$ws.Range("N85").Value = -7402

# Row 132
$ws.Range("H132").Value = 1749.3125
$ws.Range("I132").Value = 1473.6666
$ws.Range("K132").Value = 4420.9998
$ws.Range("M132").Value = -1890.9998

# Row 136
$ws.Range("H136").Value = 3073.111
$ws.Range("I136").Value = 2832.25
$ws.Range("K136").Value = 8496.75
$ws.Range("M136").Value = -5946.75
